# Updated symbol list on Mon Dec 26 14:32:05 UTC 2022 with GitHub Actions
#
# Applies the cell-level changes captured by the commit diff to the
# "cryptos" worksheet: refreshed prices in column D (kept as literal
# text so formats like trailing zeros / leading zeros survive), a
# re-ranked block of coin rows (B/C/D/E for rows 10-18), and a couple
# of ranking-label tweaks in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing a text number-format
# so that numeric-looking strings (e.g. "243.10", "0.05900") are kept
# verbatim as text instead of being coerced into floating point numbers
# (which would silently drop trailing zeros / add FP noise).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ----- Column D price refresh (value-only updates) -----
Set-TextValue $ws.Range("D2") "243.10"
Set-TextValue $ws.Range("D4") "5.409"
Set-TextValue $ws.Range("D5") "0.05901"
Set-TextValue $ws.Range("D7") "6.577"
Set-TextValue $ws.Range("D8") "0.8109"
Set-TextValue $ws.Range("D9") "0.9224"

# ----- Rows 10-18: the coin ranking shifted by one position, so every
# row's Coin / Link / Price / Volume(1h) text moved down one slot and
# a new coin ("One") was appended at the bottom of the block. -----
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1412"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D11") "0.07400"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D12") "0.03253"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03056"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09331"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Range("D15") "3.849"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D16") "0.001577"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Range("D17") "0.04674"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D18") "0.0005932"
$ws.Range("E18").Value = "17OneONE"

# ----- Remaining column D price refreshes / column E label tweaks -----
Set-TextValue $ws.Range("D19") "0.005888"

Set-TextValue $ws.Range("D20") "0.001265"
$ws.Range("E20").Value = "19BitKanKANBestin24h"

Set-TextValue $ws.Range("D21") "0.004904"
Set-TextValue $ws.Range("D22") "0.00009505"
Set-TextValue $ws.Range("D23") "3.611"
Set-TextValue $ws.Range("D25") "0.3231"

Set-TextValue $ws.Range("D40") "0.03954"
Set-TextValue $ws.Range("D41") "0.006180"
Set-TextValue $ws.Range("D42") "0.1073"
Set-TextValue $ws.Range("D43") "0.002541"
Set-TextValue $ws.Range("D44") "0.009282"
Set-TextValue $ws.Range("D45") "0.00005176"
Set-TextValue $ws.Range("D47") "0.7502"
Set-TextValue $ws.Range("D48") "0.002287"
